$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r1 = $ws.Range("C11:C16")
$r1.Style = "20% - Accent1"
$r1.Borders.Color = 8355711
$r1.Borders.Weight = 2
$r1.Borders.LineStyle = 1
$r1.Style = "20% - Accent1"
